# Applies the "text file with credentials / resolved file conflicts in
# formula data downloads" edit to slide 9 ("Download Data") of the
# presentation, and refreshes the cached date footer fields across all
# slide layouts/master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 9: rewrite the instructional bullet list in the content
#    placeholder.  The old 4 paragraphs describing how to fill in NULL
#    values / country.dir paths are replaced by a single paragraph
#    instructing the reader to locate & edit DHIS2details.txt.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(9)
$body  = $slide.Shapes.Item("Content Placeholder 2")
$tr    = $body.TextFrame.TextRange

# Paragraph 2 currently holds " At top of file (under setup), fill in NULL
# values.  "; paragraphs 3-5 hold the Example/country.dir/note text that is
# being dropped.  Replace paragraph 2's text wholesale, then delete the
# following 3 paragraphs.
$tr.Paragraphs(2, 1).Text = "Locate file DHIS2details.txt and modify as needed (then save)"
$tr.Paragraphs(3, 1).Delete()
$tr.Paragraphs(3, 1).Delete()
$tr.Paragraphs(3, 1).Delete()

# Recolor just the "DHIS2details.txt " portion of the new paragraph red.
$para2    = $tr.Paragraphs(2, 1)
$redStart = $para2.Start + 12
$redLen   = 17
$redRange = $tr.Characters($redStart, $redLen)
$redRange.Font.Color.RGB = 255

# ---------------------------------------------------------------------------
# 2. Slide 9: move the screenshot up, and remove the two small red "C"
#    right-arrow callouts that pointed at the now-removed text.
# ---------------------------------------------------------------------------
$pic = $slide.Shapes.Item("Picture 10")
$pic.Top = 226.2114960629921

$slide.Shapes.Item("Right Arrow 11").Delete()
$slide.Shapes.Item("Right Arrow 12").Delete()

# ---------------------------------------------------------------------------
# 3. Slide 9: shift the brace/label annotations for "date range of data"
#    and "download date" up to follow the repositioned screenshot.
# ---------------------------------------------------------------------------
$slide.Shapes.Item("Right Brace 14").Top = 332.7672440944882
$slide.Shapes.Item("TextBox 15").Top     = 433.22787401574806
$slide.Shapes.Item("Right Brace 16").Top = 361.83094488188976
$slide.Shapes.Item("TextBox 17").Top     = 433.1472440944882

# ---------------------------------------------------------------------------
# 4. Refresh the cached "today" date shown in the footer placeholder of
#    every slide layout and the slide master (a cosmetic field-recalc that
#    PowerPoint performs whenever the file is touched/saved again).
# ---------------------------------------------------------------------------
$newDate = "8/16/2021"

foreach ($master in $p.Designs) {
    $sm = $master.SlideMaster
    foreach ($shp in $sm.Shapes) {
        if ($shp.HasTextFrame -eq -1) {
            $ftr = $shp.TextFrame.TextRange
            if ($ftr.Text -match "/") {
                $ftr.Text = $newDate
            }
        }
    }
    for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
        $layout = $sm.CustomLayouts.Item($li)
        foreach ($shp in $layout.Shapes) {
            if ($shp.HasTextFrame -eq -1) {
                $ftr = $shp.TextFrame.TextRange
                if ($ftr.Text -match "/") {
                    $ftr.Text = $newDate
                }
            }
        }
    }
}
